# Applies the changes described in the commit:
#   "Added first task function and implemented queue. Added lots of GUI
#    elements, all working as expected for now."
#
# Concretely this adds a new "Xlarge" typography row on the Typography
# sheet, and adds/updates several rows of translated strings on the
# Translation sheet.

$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------
# Typography sheet: add a new "Xlarge" row (row 7)
# ---------------------------------------------------------------------
$wsTypo.Range("B7:J7").Style = "Normal"
$wsTypo.Range("B7").Value = "Xlarge"
$wsTypo.Range("C7").Value = "isocpeur.ttf"
$wsTypo.Range("D7").Value = 40
$wsTypo.Range("E7").Value = 4
$wsTypo.Range("F7").Value = "?"
$wsTypo.Range("I7").Value = "0-9"

# ---------------------------------------------------------------------
# Translation sheet: update existing rows
# ---------------------------------------------------------------------

# Row 11 ("Duration") GB text is unchanged, but the SI translation
# changes from "Dolzina" to "Trajanje"
$wsTrans.Range("G11").Value = "Trajanje"

# Rows 17-28: typography name changes from "Default" to "Large", and
# the SI translation ("New Text" placeholder) gets replaced with real
# translations.
$wsTrans.Range("C17").Value = "Large"
$wsTrans.Range("G17").Value = "Zamik"

$wsTrans.Range("C18").Value = "Large"
$wsTrans.Range("G18").Value = "Trajanje"

$wsTrans.Range("C19").Value = "Large"
$wsTrans.Range("G19").Value = "ms"

$wsTrans.Range("C20").Value = "Large"
$wsTrans.Range("G20").Value = "ms"

$wsTrans.Range("C21").Value = "Large"
$wsTrans.Range("G21").Value = "Zamik"

$wsTrans.Range("C22").Value = "Large"
$wsTrans.Range("G22").Value = "Trajanje"

$wsTrans.Range("C23").Value = "Large"
$wsTrans.Range("G23").Value = "ms"

$wsTrans.Range("C24").Value = "Large"
$wsTrans.Range("G24").Value = "ms"

$wsTrans.Range("C25").Value = "Large"
$wsTrans.Range("G25").Value = "Zamik"

$wsTrans.Range("C26").Value = "Large"
$wsTrans.Range("G26").Value = "Trajanje"

$wsTrans.Range("C27").Value = "Large"
$wsTrans.Range("G27").Value = "ms"

$wsTrans.Range("C28").Value = "Large"
$wsTrans.Range("G28").Value = "ms"

# ---------------------------------------------------------------------
# Translation sheet: brand new rows 29-36
# ---------------------------------------------------------------------

function Set-TranslationRow($row, $id, $typography, $alignment, $direction, $gb, $si) {
    $rng = $wsTrans.Range("B" + $row + ":G" + $row)
    $rng.Style = "Normal"
    $wsTrans.Range("B" + $row).Value = $id
    $wsTrans.Range("C" + $row).Value = $typography
    $wsTrans.Range("D" + $row).Value = $alignment
    $wsTrans.Range("E" + $row).Value = $direction
    $wsTrans.Range("F" + $row).NumberFormat = "@"
    $wsTrans.Range("F" + $row).Value = $gb
    $wsTrans.Range("F" + $row).Style = "Normal"
    $wsTrans.Range("G" + $row).NumberFormat = "@"
    $wsTrans.Range("G" + $row).Value = $si
    $wsTrans.Range("G" + $row).Style = "Normal"
}

Set-TranslationRow 29 "SingleUseId37" "Xlarge" "Center" "LTR" "<value> pcs" "<value>"
Set-TranslationRow 30 "SingleUseId38" "Xlarge" "Left"   "LTR" "0"            "-"
Set-TranslationRow 31 "SingleUseId39" "Large"  "Center" "LTR" "<value> mm"  "<value>"
Set-TranslationRow 32 "SingleUseId40" "Large"  "Left"   "LTR" "0"            "-"
Set-TranslationRow 33 "SingleUseId41" "Large"  "Center" "LTR" "<value> mm"  "<value>"
Set-TranslationRow 34 "SingleUseId42" "Large"  "Left"   "LTR" "0"            "-"
Set-TranslationRow 35 "SingleUseId43" "Large"  "Left"   "LTR" "Length setup" "Rele 3"
Set-TranslationRow 36 "SingleUseId44" "Large"  "Left"   "LTR" "mm"           "ms"

Write-Host "Edit applied"
